{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the \"x=26 \" bullet (the random-number result that needs updating) and the two\n// picture-only paragraphs directly underneath \"x=26 \" / \"y=40 \" that illustrate the old\n// random-number-generator rolls; those screenshots are being dropped from the write-up.\nlet xParagraph = null;\nlet xPictureParagraph = null;\nlet yPictureParagraph = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  if (text.indexOf(\"x=26\") !== -1) {\n    xParagraph = para;\n    xPictureParagraph = paragraphs.items[i + 1];\n  } else if (text.indexOf(\"y=40\") !== -1) {\n    yPictureParagraph = paragraphs.items[i + 1];\n  }\n}\n\n// 1) \"x=26 \" -> \"x=34 \"\nconst xSearchResults = xParagraph.search(\"26\", { matchCase: true, matchWholeWord: false });\nxSearchResults.load(\"items\");\nawait context.sync();\nxSearchResults.items[0].insertText(\"34\", \"Replace\");\n\n// 2) Drop the two screenshot-only paragraphs.\nxPictureParagraph.delete();\nyPictureParagraph.delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"x=26 \" -> \"x=34 \" in the first numbered item.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"x=26\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"x=34\"\n$find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$find.Replacement.Text, [ref]2)\n\n# 2) Remove the two paragraphs that contain only the \"True Random Number Generator\"\n#    screenshots backing up the x=26 and y=40 rolls (the screenshots are no longer needed).\n$pXPicture = $d.Paragraphs(3)\n$pYPicture = $d.Paragraphs(5)\n$pYPicture.Range.Delete()\n$pXPicture.Range.Delete()\n"}
